$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.230.75"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.14%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.596.71"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.63%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.03%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'211.83"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -0.11%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.36%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.01%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  -0.38%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.0606"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +0.40%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'18.94"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -1.38%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.0854"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.96%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'1.822.26"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +0.71%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'1.646.81"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +3.74%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'4.01"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -0.09%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('E15').Value = "'  -2.20%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D17').Value = "'26.217.74"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.10%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'227.80"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +6.51%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'  -0.71%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'7.57"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +3.72%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'  +0.00%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'4.24"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E23').Value = "'  +0.50%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = "'  -0.69%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'145.50"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +1.14%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D27').Value = "'6.96"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -0.56%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  +0.89%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  +1.64%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'0.0492"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -0.69%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  +0.20%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  +0.28%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'1.446.30"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +3.67%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'2.94"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +0.00%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  -0.28%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'  +0.49%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'0.563"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -3.72%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  -1.57%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'0.817"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -0.18%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'5.73"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -2.29%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  +0.03%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'2.18"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +1.94%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'  -0.88%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'1.734.12"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +0.67%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'  -0.99%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'60.40"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -0.92%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'87.57"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +1.84%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  -0.63%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('E50').Value = "'  -0.11%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'7.41"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +0.67%  "
$ws.Range('E51').Style = 'Normal'
